# Updated cryptos list on Thu Jun 13 10:29:37 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "67.373.73"
$ws.Range("E2").Value = "  -0.70%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.483.07"
$ws.Range("E3").Value = "  -1.78%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.09%  "

# Row 5 - BNB
$ws.Range("D5").Value = "604.49"
$ws.Range("E5").Value = "  -1.96%  "

# Row 6 - Solana
$ws.Range("D6").Value = "150.77"
$ws.Range("E6").Value = "  -1.23%  "

# Row 7 - LidoStakedEther
$ws.Range("D7").Value = "3.481.30"
$ws.Range("E7").Value = "  -1.76%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.03%  "

# Row 9 - XRP
$ws.Range("D9").Value = "0.486"
$ws.Range("E9").Value = "  +0.79%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +2.43%  "

# Row 11 - Toncoin
$ws.Range("D11").Value = "7.54"
$ws.Range("E11").Value = "  +6.56%  "

# Row 12 - Cardano
$ws.Range("D12").Value = "0.431"
$ws.Range("E12").Value = "  +0.93%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  -1.82%  "

# Row 14 - Avalanche
$ws.Range("D14").Value = "32.02"
$ws.Range("E14").Value = "  -0.58%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "4.067.53"
$ws.Range("E15").Value = "  -1.93%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "3.475.81"
$ws.Range("E16").Value = "  -2.08%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "67.357.12"
$ws.Range("E17").Value = "  -0.44%  "

# Row 18 - TRON
$ws.Range("E18").Value = "  -0.08%  "

# Row 19 - Polkadot
$ws.Range("D19").Value = "6.47"
$ws.Range("E19").Value = "  +0.79%  "

# Row 20 - Chainlink
$ws.Range("D20").Value = "15.41"
$ws.Range("E20").Value = "  +0.18%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  +1.86%  "

# Row 22 - BitcoinCash
$ws.Range("D22").Value = "445.88"
$ws.Range("E22").Value = "  -0.48%  "

# Row 23 - Polygon
$ws.Range("D23").Value = "0.626"
$ws.Range("E23").Value = "  +0.01%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "78.91"
$ws.Range("E24").Value = "  +1.78%  "

# Row 25 - now Dai (was WrappedeETH)
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  -0.09%  "

# Row 26 - now WrappedeETH (was Dai)
$ws.Range("B26").Value = "WrappedeETH"
$ws.Range("C26").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D26").Value = "3.623.99"
$ws.Range("E26").Value = "  -1.78%  "

# Row 27 - PEPE
$ws.Range("E27").Value = "  -5.00%  "

# Row 28 - RenderToken
$ws.Range("D28").Value = "8.65"
$ws.Range("E28").Value = "  +0.00%  "

# Row 29 - InternetComputer(DFINITY)
$ws.Range("D29").Value = "9.94"
$ws.Range("E29").Value = "  -3.69%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  -1.41%  "

# Row 31 - Fetch.AI
$ws.Range("E31").Value = "  +2.47%  "

# Row 32 - Kaspa
$ws.Range("D32").Value = "0.170"
$ws.Range("E32").Value = "  +1.27%  "

# Row 33 - Binance-PegBSC-USD
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  -0.17%  "

# Row 34 - EthereumClassic
$ws.Range("D34").Value = "25.57"
$ws.Range("E34").Value = "  -1.69%  "

# Row 35 - NEARProtocol
$ws.Range("D35").Value = "6.12"
$ws.Range("E35").Value = "  -1.64%  "

# Row 36 - ImmutableX
$ws.Range("D36").Value = "1.85"
$ws.Range("E36").Value = "  -0.45%  "

# Row 37 - RenzoRestakedETH
$ws.Range("D37").Value = "3.475.69"
$ws.Range("E37").Value = "  -1.76%  "

# Row 38 - Aptos
$ws.Range("E38").Value = "  -1.04%  "

# Row 39 - USDe
$ws.Range("E39").Value = "  +0.01%  "

# Row 40 - Stacks
$ws.Range("E40").Value = "  +5.13%  "

# Row 41 - Monero
$ws.Range("D41").Value = "177.50"
$ws.Range("E41").Value = "  +0.30%  "

# Row 42 - FirstDigitalUSD
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  -0.16%  "

# Row 43 - Hedera
$ws.Range("D43").Value = "0.0895"
$ws.Range("E43").Value = "  -0.40%  "

# Row 44 - Filecoin
$ws.Range("E44").Value = "  -0.33%  "

# Row 45 - Mantle
$ws.Range("E45").Value = "  +0.17%  "

# Row 46 - InjectiveProtocol
$ws.Range("D46").Value = "30.01"
$ws.Range("E46").Value = "  +4.71%  "

# Row 47 - OKB
$ws.Range("D47").Value = "46.49"
$ws.Range("E47").Value = "  +2.38%  "

# Row 48 - ONDO
$ws.Range("E48").Value = "  +0.37%  "

# Row 49 - dogwifhat
$ws.Range("D49").Value = "2.53"
$ws.Range("E49").Value = "  -5.33%  "

# Row 50 - Cosmos
$ws.Range("D50").Value = "7.59"
$ws.Range("E50").Value = "  -0.45%  "

# Row 51 - TheGraph
$ws.Range("D51").Value = "0.251"
$ws.Range("E51").Value = "  -0.45%  "
